$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-07 01:41:02"

for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
